# Auto-generated edit script: apply updated market-price figures to Goblin Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value2 = 45793.43
$ws.Range("J3").Value2 = 45793.43
$ws.Range("L3").Value2 = 45793.43
$ws.Range("N3").Value2 = -46021.43
$ws.Range("H20").Value2 = 3124
$ws.Range("I20").Value2 = 3124
$ws.Range("K20").Value2 = 3124
$ws.Range("M20").Value2 = -2894
$ws.Range("H35").Value2 = 3124
$ws.Range("I35").Value2 = 3124
$ws.Range("K35").Value2 = 3124
$ws.Range("M35").Value2 = -2745
$ws.Range("H40").Value2 = 3999.524
$ws.Range("I40").Value2 = 1841.3334
$ws.Range("J40").Value2 = 4862.8
$ws.Range("K40").Value2 = 1841.3334
$ws.Range("L40").Value2 = 4862.8
$ws.Range("M40").Value2 = -1666.3334
$ws.Range("N40").Value2 = -5212.8
$ws.Range("H102").Value2 = 45793.43
$ws.Range("J102").Value2 = 45793.43
$ws.Range("L102").Value2 = 45793.43
$ws.Range("N102").Value2 = -52283.43
$ws.Range("H108").Value2 = 75000
$ws.Range("J108").Value2 = 75000
$ws.Range("L108").Value2 = 75000
$ws.Range("N108").Value2 = -82680
$ws.Range("H116").Value2 = 6500.6
$ws.Range("J116").Value2 = 11501.5
$ws.Range("L116").Value2 = 11501.5
$ws.Range("N116").Value2 = -18385.5
$ws.Range("H133").Value2 = 104997.5
$ws.Range("J133").Value2 = 104997.5
$ws.Range("L133").Value2 = 104997.5
$ws.Range("N133").Value2 = -115117.5
$ws.Range("H136").Value2 = 116923
$ws.Range("J136").Value2 = 116923
$ws.Range("L136").Value2 = 116923
$ws.Range("N136").Value2 = -127123
$ws.Range("H137").Value2 = 2199.2632
$ws.Range("J137").Value2 = 2058.1667
$ws.Range("L137").Value2 = 6174.500100000001
$ws.Range("N137").Value2 = -11274.5001
$ws.Range("H138").Value2 = 1381208.8
$ws.Range("I138").Value2 = 10837.818
$ws.Range("J138").Value2 = 1596552.6
$ws.Range("K138").Value2 = 32513.454
$ws.Range("L138").Value2 = 4789657.800000001
$ws.Range("M138").Value2 = -27373.454
$ws.Range("N138").Value2 = -4799937.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value2 = 60000
$ws.Range("J113").Value2 = 60000
$ws.Range("L113").Value2 = 60000
$ws.Range("N113").Value2 = -68678
$ws.Range("H133").Value2 = 35994
$ws.Range("J133").Value2 = 36192.8
$ws.Range("L133").Value2 = 36192.8
$ws.Range("N133").Value2 = -41252.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 3339.4707
$ws.Range("I99").Value2 = 2160.4546
$ws.Range("K99").Value2 = 2160.4546
$ws.Range("M99").Value2 = -662.4546
$ws.Range("H105").Value2 = 10014.77
$ws.Range("I105").Value2 = 14025.875
$ws.Range("J105").Value2 = 3597
$ws.Range("K105").Value2 = 14025.875
$ws.Range("L105").Value2 = 3597
$ws.Range("M105").Value2 = -12278.875
$ws.Range("N105").Value2 = -7091
$ws.Range("H107").Value2 = 5699.2856
$ws.Range("I107").Value2 = 2948
$ws.Range("K107").Value2 = 2948
$ws.Range("M107").Value2 = -1028
$ws.Range("H132").Value2 = 89904.836
$ws.Range("J132").Value2 = 89904.836
$ws.Range("L132").Value2 = 89904.836
$ws.Range("N132").Value2 = -100024.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value2 = 1537.0769
$ws.Range("I10").Value2 = 1239.5
$ws.Range("J10").Value2 = 2529
$ws.Range("K10").Value2 = 1239.5
$ws.Range("L10").Value2 = 2529
$ws.Range("M10").Value2 = -1100.5
$ws.Range("N10").Value2 = -2807
$ws.Range("H31").Value2 = 4931.6665
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 4931.6665
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 4931.6665
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value2 = -5521.6665
$ws.Range("H32").Value2 = 3245.5
$ws.Range("I32").Value2 = 4294
$ws.Range("K32").Value2 = 4294
$ws.Range("M32").Value2 = -3978
$ws.Range("H34").Value2 = 4931.6665
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 4931.6665
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 4931.6665
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value2 = -5335.6665
$ws.Range("H60").Value2 = 29843.4
$ws.Range("J60").Value2 = 52108.5
$ws.Range("L60").Value2 = 52108.5
$ws.Range("N60").Value2 = -53130.5
$ws.Range("H103").Value2 = 14604.8
$ws.Range("I103").Value2 = 14604.8
$ws.Range("K103").Value2 = 14604.8
$ws.Range("M103").Value2 = -13432.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value2 = 413.2
$ws.Range("I33").Value2 = 346.33334
$ws.Range("J33").Value2 = 429.91666
$ws.Range("K33").Value2 = 2078.00004
$ws.Range("L33").Value2 = 2579.49996
$ws.Range("M33").Value2 = -1795.00004
$ws.Range("N33").Value2 = -3145.49996
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value2 = 2300
$ws.Range("I97").Value2 = 1400
$ws.Range("J97").Value2 = 5000
$ws.Range("K97").Value2 = 4200
$ws.Range("L97").Value2 = 15000
$ws.Range("M97").Value2 = -3704
$ws.Range("N97").Value2 = -15992
$ws.Range("H112").Value2 = 6415
$ws.Range("I112").Value2 = 4215
$ws.Range("J112").Value2 = 7515
$ws.Range("K112").Value2 = 12645
$ws.Range("L112").Value2 = 22545
$ws.Range("M112").Value2 = -11537
$ws.Range("N112").Value2 = -24761
$ws.Range("H139").Value2 = 4917.773
$ws.Range("I139").Value2 = 4835.5454
$ws.Range("K139").Value2 = 14506.6362
$ws.Range("M139").Value2 = -9366.636200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 20849382
$ws.Range("I70").Value2 = 33350732
$ws.Range("J70").Value2 = 13798
$ws.Range("K70").Value2 = 33350732
$ws.Range("L70").Value2 = 13798
$ws.Range("M70").Value2 = -33350462
$ws.Range("N70").Value2 = -14338
$ws.Range("H73").Value2 = 20849382
$ws.Range("I73").Value2 = 33350732
$ws.Range("J73").Value2 = 13798
$ws.Range("K73").Value2 = 33350732
$ws.Range("L73").Value2 = 13798
$ws.Range("M73").Value2 = -33349796
$ws.Range("N73").Value2 = -15670
$ws.Range("H80").Value2 = 5685.122
$ws.Range("I80").Value2 = 3760.7896
$ws.Range("J80").Value2 = 7347.0454
$ws.Range("K80").Value2 = 3760.7896
$ws.Range("L80").Value2 = 7347.0454
$ws.Range("M80").Value2 = -2762.7896
$ws.Range("N80").Value2 = -9343.045399999999
$ws.Range("H83").Value2 = 5685.122
$ws.Range("I83").Value2 = 3760.7896
$ws.Range("J83").Value2 = 7347.0454
$ws.Range("K83").Value2 = 18803.948
$ws.Range("L83").Value2 = 36735.227
$ws.Range("M83").Value2 = -13811.948
$ws.Range("N83").Value2 = -46719.227
$ws.Range("H107").Value2 = 1523.6666
$ws.Range("I107").Value2 = 425.5
$ws.Range("J107").Value2 = 2199.4614
$ws.Range("K107").Value2 = 425.5
$ws.Range("L107").Value2 = 2199.4614
$ws.Range("M107").Value2 = 1494.5
$ws.Range("N107").Value2 = -6039.4614
$ws.Range("H117").Value2 = 40000
$ws.Range("J117").Value2 = 40000
$ws.Range("L117").Value2 = 40000
$ws.Range("N117").Value2 = -46884
$ws.Range("H123").Value2 = 38299.332
$ws.Range("J123").Value2 = 38299.332
$ws.Range("L123").Value2 = 38299.332
$ws.Range("N123").Value2 = -43199.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value2 = 4686.6
$ws.Range("I32").Value2 = 4686.6
$ws.Range("K32").Value2 = 4686.6
$ws.Range("M32").Value2 = -4369.6
$ws.Range("H68").Value2 = 7019.55
$ws.Range("J68").Value2 = 8308.333000000001
$ws.Range("L68").Value2 = 8308.333000000001
$ws.Range("N68").Value2 = -9806.333000000001
$ws.Range("H71").Value2 = 7019.55
$ws.Range("J71").Value2 = 8308.333000000001
$ws.Range("L71").Value2 = 41541.665
$ws.Range("N71").Value2 = -49029.665
$ws.Range("H98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("N98").ClearContents()
$ws.Range("H103").Value2 = 13908.167
$ws.Range("J103").Value2 = 13908.167
$ws.Range("L103").Value2 = 13908.167
$ws.Range("N103").Value2 = -16252.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value2 = 46665
$ws.Range("I99").Value2 = 40000
$ws.Range("K99").Value2 = 40000
$ws.Range("M99").Value2 = -37005
$ws.Range("H122").Value2 = 3633.9556
$ws.Range("I122").Value2 = 2557.6428
$ws.Range("J122").Value2 = 4120.032
$ws.Range("K122").Value2 = 7672.928400000001
$ws.Range("L122").Value2 = 12360.096
$ws.Range("M122").Value2 = -5222.928400000001
$ws.Range("N122").Value2 = -17260.096
$ws.Range("H126").Value2 = 2507.2222
$ws.Range("I126").Value2 = 1959.7826
$ws.Range("K126").Value2 = 5879.3478
$ws.Range("M126").Value2 = -3409.3478
$ws.Range("H133").Value2 = 83499.5
$ws.Range("J133").Value2 = 83499.5
$ws.Range("L133").Value2 = 83499.5
$ws.Range("N133").Value2 = -93619.5
